# Update the East-Asian and Complex-Script font fallbacks used across the
# document's paragraph styles (mirrors a docx/html regeneration under
# vignettes: the East Asian fallback moves from "DejaVu Sans" to "Tahoma",
# and the Complex-Script fallback ("DejaVu Sans") becomes explicit on a few
# styles that previously inherited it).

$d = $word.ActiveDocument

# Styles whose East Asian fallback font changes from "DejaVu Sans" to "Tahoma".
$eastAsianStyles = @("Normal", "Heading")
foreach ($styleName in $eastAsianStyles) {
    $style = $d.Styles($styleName)
    $style.Font.NameFarEast = "Tahoma"
}

# Styles that gain an explicit Complex-Script font (same value they already
# inherited: "DejaVu Sans") in their own run properties.
$csStyles = @("List", "Caption", "Index")
foreach ($styleName in $csStyles) {
    $style = $d.Styles($styleName)
    $style.Font.NameBi = "DejaVu Sans"
}
